# Add new columns I (I0) and J (IF) to the sheet, mirroring the diff:
# - Header row 1: I1 = "I0", J1 = "IF"
# - Rows 2-24: numeric values for I and J

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, border, centered) from an existing header cell (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for column I and J, rows 2..24
$values = @(
    @(6, 6),    # row 2
    @(8, 8),    # row 3
    @(10, 10),  # row 4
    @(10, 10),  # row 5
    @(7, 8),    # row 6
    @(8, 8),    # row 7
    @(6, 7),    # row 8
    @(5, 5),    # row 9
    @(8, 8),    # row 10
    @(6, 6),    # row 11
    @(8, 9),    # row 12
    @(9, 9),    # row 13
    @(7, 7),    # row 14
    @(9, 9),    # row 15
    @(3, 3),    # row 16
    @(6, 6),    # row 17
    @(8, 8),    # row 18
    @(8, 9),    # row 19
    @(5, 5),    # row 20
    @(9, 9),    # row 21
    @(4, 4),    # row 22
    @(7, 7),    # row 23
    @(4, 4)     # row 24
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
